# Excel COM-interop script implementing the commit:
#   Insert a new "債務" (debt) sheet where "事業投資" used to sit, with a
#   fresh (header-only) layout, and move the old "事業投資" data to a brand
#   new sheet appended at the end of the workbook.

$wb = $excel.ActiveWorkbook

# 1) Grab the sheet that is about to become "債務" and remember its old
#    (事業投資) data before we touch anything.
$invest = $wb.Worksheets.Item("事業投資")

$ownerVal    = $invest.Range("B2").Value
$companyVal  = $invest.Range("C2").Value
$addressVal  = $invest.Range("D2").Value
$totalVal    = $invest.Range("E2").Value
$regDateVal  = $invest.Range("F2").Value
$regReasonVal= $invest.Range("G2").Value
$propCatVal  = $invest.Range("H2").Value
$categoryVal = $invest.Range("I2").Value
$dateVal     = $invest.Range("J2").Value
$legNameVal  = $invest.Range("K2").Value
$legIdVal    = $invest.Range("L2").Value
$srcFileVal  = $invest.Range("M2").Value
$indexVal    = $invest.Range("N2").Value

# 2) Rename the existing sheet to "債務" (keeps its sheetId/rId slot) and
#    wipe its contents, then lay down the new header-only row.
$debt = $invest
$debt.Cells.Clear()
$debt.Name = "債務"

$debt.Range("B1").Value = "species"
$debt.Range("C1").Value = "debtor"
$debt.Range("D1").Value = "owner"
$debt.Range("E1").Value = "total"
$debt.Range("F1").Value = "register_date"
$debt.Range("G1").Value = "register_reason"
$debt.Range("H1").Value = "property_category"
$debt.Range("I1").Value = "category"
$debt.Range("J1").Value = "date"
$debt.Range("K1").Value = "legislator_name"
$debt.Range("L1").Value = "legislator_id"
$debt.Range("M1").Value = "source_file"
$debt.Range("N1").Value = "index"

$headerRow = $debt.Range("B1:N1")
$headerRow.Font.Bold = $true
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160
$headerRow.Borders.LineStyle = 1

# 3) Append a brand new sheet (gets the next sheetId/rId) named "事業投資"
#    at the end of the workbook, carrying the data that used to live in
#    the now-renamed "債務" sheet.
$newInvest = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newInvest.Name = "事業投資"

$newInvest.Range("B1").Value = "owner"
$newInvest.Range("C1").Value = "company"
$newInvest.Range("D1").Value = "address"
$newInvest.Range("E1").Value = "total"
$newInvest.Range("F1").Value = "register_date"
$newInvest.Range("G1").Value = "register_reason"
$newInvest.Range("H1").Value = "property_category"
$newInvest.Range("I1").Value = "category"
$newInvest.Range("J1").Value = "date"
$newInvest.Range("K1").Value = "legislator_name"
$newInvest.Range("L1").Value = "legislator_id"
$newInvest.Range("M1").Value = "source_file"
$newInvest.Range("N1").Value = "index"

$newHeaderRow = $newInvest.Range("B1:N1")
$newHeaderRow.Font.Bold = $true
$newHeaderRow.HorizontalAlignment = -4108
$newHeaderRow.VerticalAlignment = -4160
$newHeaderRow.Borders.LineStyle = 1

$newInvest.Range("A2").Value = $indexVal
$newInvest.Range("B2").Value = $ownerVal
$newInvest.Range("C2").Value = $companyVal
$newInvest.Range("D2").Value = $addressVal
$newInvest.Range("E2").Value = $totalVal
$newInvest.Range("F2").Value = $regDateVal
$newInvest.Range("G2").Value = $regReasonVal
$newInvest.Range("H2").Value = $propCatVal
$newInvest.Range("I2").Value = $categoryVal
$newInvest.Range("J2").Value = $dateVal
$newInvest.Range("K2").Value = $legNameVal
$newInvest.Range("L2").Value = $legIdVal
$newInvest.Range("M2").Value = $srcFileVal
$newInvest.Range("N2").Value = $indexVal

$newInvest.Range("A2").Font.Bold = $true
$newInvest.Range("A2").HorizontalAlignment = -4108
$newInvest.Range("A2").VerticalAlignment = -4160
$newInvest.Range("A2").Borders.LineStyle = 1
